# Extend the "SanCesarioSP" daily-counts sheet with new rows 302-328
# (data updated through 25 July 2021 / commit: "aggiornamento fino a 28 luglio").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data: row number, date serial (col A), nuovi pos. (col B),
# somma mobile 7gg. (col C), somma mobile 7gg. per 100mila abitanti (col D)
$newRows = @(
    @(302, 44376, 0, 1, 15.22997258604935),
    @(303, 44377, 0, 1, 15.22997258604935),
    @(304, 44378, 0, 0, 0),
    @(305, 44379, 0, 0, 0),
    @(306, 44380, 0, 0, 0),
    @(307, 44381, 0, 0, 0),
    @(308, 44382, 0, 0, 0),
    @(309, 44383, 0, 0, 0),
    @(310, 44384, 0, 0, 0),
    @(311, 44385, 0, 0, 0),
    @(312, 44386, 0, 0, 0),
    @(313, 44387, 0, 0, 0),
    @(314, 44388, 0, 0, 0),
    @(315, 44389, 0, 0, 0),
    @(316, 44390, 0, 0, 0),
    @(317, 44391, 0, 0, 0),
    @(318, 44392, 0, 0, 0),
    @(319, 44393, 0, 0, 0),
    @(320, 44394, 0, 0, 0),
    @(321, 44395, 0, 0, 0),
    @(322, 44396, 0, 0, 0),
    @(323, 44397, 0, 0, 0),
    @(324, 44398, 0, 0, 0),
    @(325, 44399, 1, 1, 15.22997258604935),
    @(326, 44400, 0, 1, 15.22997258604935),
    @(327, 44401, 0, 1, 15.22997258604935),
    @(328, 44402, 0, 1, 15.22997258604935)
)

foreach ($r in $newRows) {
    $rowNum = $r[0]

    # Column A carries the date style (s="2") used by every row above it;
    # copy that formatting from the prior row before writing the new values.
    $ws.Range("A" + ($rowNum - 1)).Copy()
    $ws.Range("A" + $rowNum).PasteSpecial(-4122)

    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
}

$excel.CutCopyMode = $false
